$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range('D2').Value = '61.655.91'
$ws.Range('E2').Value = '  -3.23%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range('D3').Value = '2.482.27'
$ws.Range('E3').Value = '  -5.64%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range('E4').Value = '  +0.06%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range('D5').Value = '553.41'
$ws.Range('E5').Value = '  -4.45%  '

# Row 6: 'Solana' -> 'Solana'
$ws.Range('D6').Value = '146.98'
$ws.Range('E6').Value = '  -5.25%  '

# Row 7: 'USDC' -> 'USDC'
$ws.Range('E7').Value = '  +0.02%  '

# Row 8: 'XRP' -> 'XRP'
$ws.Range('D8').Value = '0.601'
$ws.Range('E8').Value = '  -3.47%  '

# Row 9: 'LidoStakedEther' -> 'LidoStakedEther'
$ws.Range('D9').Value = '2.481.01'
$ws.Range('E9').Value = '  -5.55%  '

# Row 10: 'Dogecoin' -> 'Dogecoin'
$ws.Range('E10').Value = '  -8.29%  '

# Row 11: 'Toncoin' -> 'TRON'
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '0.154'
$ws.Range('E11').Value = '  -1.38%  '

# Row 12: 'TRON' -> 'Toncoin'
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').Value = '5.43'
$ws.Range('E12').Value = '  -6.34%  '

# Row 13: 'Cardano' -> 'Cardano'
$ws.Range('D13').Value = '0.359'
$ws.Range('E13').Value = '  -5.76%  '

# Row 14: 'Avalanche' -> 'Avalanche'
$ws.Range('D14').Value = '26.29'
$ws.Range('E14').Value = '  -7.29%  '

# Row 15: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range('D15').Value = '2.930.76'
$ws.Range('E15').Value = '  -5.58%  '

# Row 16: 'ShibaInu' -> 'ShibaInu'
$ws.Range('D16').Value = '0.0000167'
$ws.Range('E16').Value = '  -8.44%  '

# Row 17: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range('D17').Value = '61.580.96'
$ws.Range('E17').Value = '  -3.26%  '

# Row 18: 'WrappedEther' -> 'WrappedEther'
$ws.Range('D18').Value = '2.484.17'
$ws.Range('E18').Value = '  -5.31%  '

# Row 19: 'Chainlink' -> 'Chainlink'
$ws.Range('D19').Value = '11.16'
$ws.Range('E19').Value = '  -7.77%  '

# Row 20: 'Uniswap' -> 'Uniswap'
$ws.Range('D20').Value = '7.00'
$ws.Range('E20').Value = '  -8.33%  '

# Row 21: 'Polkadot' -> 'Polkadot'
$ws.Range('D21').Value = '4.21'
$ws.Range('E21').Value = '  -6.89%  '

# Row 22: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range('D22').Value = '322.66'
$ws.Range('E22').Value = '  -6.31%  '

# Row 23: 'Dai' -> 'Dai'
$ws.Range('E23').Value = '  -0.03%  '

# Row 24: 'SuiNetwork' -> 'SuiNetwork'
$ws.Range('E24').Value = '  -5.43%  '

# Row 25: 'Litecoin' -> 'Litecoin'
$ws.Range('D25').Value = '64.14'
$ws.Range('E25').Value = '  -5.58%  '

# Row 26: 'PEPE' -> 'PEPE'
$ws.Range('D26').Value = '0.0₃0988'
$ws.Range('E26').Value = '  -8.83%  '

# Row 27: 'Fetch.AI' -> 'WrappedeETH'
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.612.36'
$ws.Range('E27').Value = '  -5.18%  '

# Row 28: 'WrappedeETH' -> 'Fetch.AI'
$ws.Range('B28').Value = 'Fetch.AI'
$ws.Range('C28').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D28').Value = '1.54'
$ws.Range('E28').Value = '  -4.07%  '

# Row 29: 'Binance-PegBSC-USD' -> 'Bittensor'
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').Value = '540.28'
$ws.Range('E29').Value = '  -10.25%  '

# Row 30: 'Bittensor' -> 'Binance-PegBSC-USD'
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.01%  '

# Row 31: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range('D31').Value = '8.39'
$ws.Range('E31').Value = '  -9.32%  '

# Row 32: 'Aptos' -> 'Aptos'
$ws.Range('D32').Value = '7.56'
$ws.Range('E32').Value = '  -6.17%  '

# Row 33: 'Kaspa' -> 'Kaspa'
$ws.Range('E33').Value = '  -5.77%  '

# Row 34: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range('D34').Value = '1.91'
$ws.Range('E34').Value = '  -7.57%  '

# Row 35: 'ImmutableX' -> 'ImmutableX'
$ws.Range('D35').Value = '1.60'
$ws.Range('E35').Value = '  -8.82%  '

# Row 36: 'RenderToken' -> 'RenderToken'
$ws.Range('D36').Value = '5.95'
$ws.Range('E36').Value = '  -9.98%  '

# Row 37: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range('D37').Value = '4.89'
$ws.Range('E37').Value = '  -9.98%  '

# Row 38: 'FirstDigitalUSD' -> 'FirstDigitalUSD'
$ws.Range('E38').Value = '  +0.11%  '

# Row 39: 'PolygonEcosystemToken' -> 'PolygonEcosystemToken'
$ws.Range('D39').Value = '0.383'
$ws.Range('E39').Value = '  -4.81%  '

# Row 40: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range('D40').Value = '18.56'
$ws.Range('E40').Value = '  -5.92%  '

# Row 41: 'Monero' -> 'Monero'
$ws.Range('D41').Value = '148.70'
$ws.Range('E41').Value = '  -0.77%  '

# Row 42: 'Stacks' -> 'Stacks'
$ws.Range('D42').Value = '1.74'
$ws.Range('E42').Value = '  -8.24%  '

# Row 43: 'USDe' -> 'USDe'
$ws.Range('E43').Value = '  +0.05%  '

# Row 44: 'OKB' -> 'OKB'
$ws.Range('D44').Value = '40.37'
$ws.Range('E44').Value = '  -3.60%  '

# Row 45: 'dogwifhat' -> 'dogwifhat'
$ws.Range('D45').Value = '2.34'
$ws.Range('E45').Value = '  -7.66%  '

# Row 46: 'Aave' -> 'Aave'
$ws.Range('D46').Value = '148.28'
$ws.Range('E46').Value = '  -6.87%  '

# Row 47: 'Filecoin' -> 'Filecoin'
$ws.Range('D47').Value = '3.63'
$ws.Range('E47').Value = '  -7.09%  '

# Row 48: 'InjectiveProtocol' -> 'InjectiveProtocol'
$ws.Range('D48').Value = '20.98'
$ws.Range('E48').Value = '  -14.84%  '

# Row 49: 'Hedera' -> 'Hedera'
$ws.Range('D49').Value = '0.0538'
$ws.Range('E49').Value = '  -8.38%  '

# Row 50: 'Mantle' -> 'Mantle'
$ws.Range('D50').Value = '0.598'
$ws.Range('E50').Value = '  -5.20%  '

# Row 51: 'Stellar' -> 'Stellar'
$ws.Range('D51').Value = '0.0949'
$ws.Range('E51').Value = '  -4.84%  '
